# feat: add 2022-Q1 data
#
# The workbook currently has 4 sheets: 2021-Q2, 2021-Q3, 2021-Q4, 总计 (Total).
# This script:
#   1. Renames the current "总计" sheet to "2022-Q1" and replaces its content
#      with the 2022-Q1 fund-holding detail rows (same shape as the other
#      2021-Qx detail sheets).
#   2. Creates a fresh copy of that sheet, renamed to "总计", and rewrites its
#      content to be the quarterly roll-up table with a new leading row for
#      2022-Q1 (the pre-existing rows shift down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: the old "总计" sheet becomes "2022-Q1"
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Make a fresh copy of it (placed right after) before overwriting its
# content -- the copy becomes the new "总计" sheet, inheriting the same
# sheet formatting (outline props, page margins, base col width, etc.)
$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item("2022-Q1 (2)")
$total.Name = "总计"

# A throw-away worksheet, added at the very end (after every real sheet) so
# it doesn't disturb sheet ordering/ids, used purely as a relay: writing a
# string-literal formula into it and then Copy / PasteSpecial-values into
# the real destination is the only way to land a cell that stores TEXT even
# when the text looks numeric ("22.84", "506005", ...) without leaving a
# quotePrefix style behind (which is what happens if you just assign
# .Value = "'22.84").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$scratch = $wb.Worksheets.Add($null, $lastSheet)
$scratchCell = $scratch.Cells.Item(1, 1)

function Set-TextValue($cell, [string]$value) {
    $escaped = $value -replace '"', '""'
    $scratchCell.Formula = '="' + $escaped + '"'
    $scratch.Range("A1").Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# ---------------------------------------------------------------------
# Step 2: rewrite "2022-Q1" with the fund holdings detail for the quarter
# ---------------------------------------------------------------------
$q1.Range("A1:D4").ClearContents()

# Extend the header/index styling (s="2") that already exists on B1:D1 and
# A2 out to the newly-used E1:H1 and A5 cells.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)  # xlPasteFormats
$q1.Range("A2").Copy()
$q1.Range("A5").PasteSpecial(-4122)  # xlPasteFormats

$q1Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $q1Headers.Length; $c++) {
    Set-TextValue $q1.Cells.Item(1, 2 + $c) $q1Headers[$c]
}

$q1Rows = @(
    @("506005", "博时科创板三年定期开放混合", "22.84", "96.44", "3.60", "0.8222", 8),
    @("011851", "天弘先进制造混合型证券投资基金A", "2.72", "91.41", "7.04", "0.1915", 3),
    @("004784", "招商稳健优选股票", "4.90", "85.49", "3.19", "0.1563", 6),
    @("011852", "天弘先进制造混合型证券投资基金C", "0.70", "91.41", "7.04", "0.0493", 3)
)

for ($r = 0; $r -lt $q1Rows.Length; $r++) {
    $row = 2 + $r
    $q1.Cells.Item($row, 1).Value = $r
    for ($c = 0; $c -lt 6; $c++) {
        Set-TextValue $q1.Cells.Item($row, 2 + $c) $q1Rows[$r][$c]
    }
    $q1.Cells.Item($row, 8).Value = $q1Rows[$r][6]
}

# ---------------------------------------------------------------------
# Step 3: rewrite "总计" with the updated roll-up (new 2022-Q1 row on top)
# ---------------------------------------------------------------------
$total.Range("A1:D4").ClearContents()

# The roll-up gains a 5th row (one new quarter), so extend the index styling
# (s="2") from A2 down to the newly-used A5 cell, same as step 2 above.
$total.Range("A2").Copy()
$total.Range("A5").PasteSpecial(-4122)  # xlPasteFormats

Set-TextValue $total.Cells.Item(1, 2) "日期"
Set-TextValue $total.Cells.Item(1, 3) "持有数量(只)"
Set-TextValue $total.Cells.Item(1, 4) "持有市值(亿元)"

$totalRows = @(
    @("2022-Q1", 4, 1.22),
    @("2021-Q4", 1, 0),
    @("2021-Q3", 3, 1.77),
    @("2021-Q2", 1, 0.07000000000000001)
)

for ($r = 0; $r -lt $totalRows.Length; $r++) {
    $row = 2 + $r
    $total.Cells.Item($row, 1).Value = $r
    Set-TextValue $total.Cells.Item($row, 2) $totalRows[$r][0]
    $total.Cells.Item($row, 3).Value = $totalRows[$r][1]
    $total.Cells.Item($row, 4).Value = $totalRows[$r][2]
}

# ---------------------------------------------------------------------
# Clean up the relay sheet
# ---------------------------------------------------------------------
$scratch.Delete() | Out-Null
